# Update Gruppe C und E
# - Sprint-Backlog!D6: re-estimate "240min" -> "120min"
# - Sprint-Backlog!E4: mark done ("Done")
# - Sprint-Backlog row 8: new backlog item (TicTacToe rules implementation), 120min
# - Selection/view state: Backlog -> row 7 selected; Sprint-Backlog -> D8 selected

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Backlog")
$ws2 = $wb.Worksheets.Item("Sprint-Backlog")

# Re-estimate the TicTacToe-Darstellung task from 240min to 120min.
$ws2.Range("D6").Value = "120min"

# Mark the "Clonen des Repositories..." task as Done.
$ws2.Range("E4").Value = "Done"

# New Sprint-Backlog row: "Implementierung der TicTacToe-Regeln" task.
$ws2.Range("A8").Value = 2
$ws2.Range("B8").Value = "Implementierung der TicTacToe-Regeln"
$ws2.Range("C8").Value = "Umsetzung der Rules im Framework. Standard TicTacToe. Mit Ausblick aus Zeit-Limit."
$ws2.Range("D8").Value = "120min"

# Restore view/selection state as left by the editor.
$ws1.Activate()
$ws1.Rows(7).Select()

$ws2.Activate()
$ws2.Range("D8").Select()
